# Updated notebook, reran simulation
# - added two new materials: "Holden" and "Rizzie Spiral" (inserted right
#   after "Spiral5")
# - renamed "Thomas Hex" -> "Matthies Hex"
# - reran the simulation, which re-emits the result table and appends two
#   more rows ("Michael-CCHex" / "Michael-SNHex") at the bottom

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full, final ordered list of material names for column B, rows 3..31
# (A holds the corresponding 1-based-minus-1 index, i.e. row-2).
$materials = @(
    "Spiral5",
    "Holden",
    "Rizzie Spiral",
    "RotRing OmegaMax-90",
    "Equal Angle",
    "Tilt Rotate",
    "CLR",
    "Rizzie Hex",
    "Matthies Hex",
    "Tilt Rotate_Partial",
    "RotRing OmegaMax-60",
    "Equal Angle_Partial",
    "Rizzie Hex_Partial",
    "ND Single",
    "RD Single",
    "TD Single",
    "Morris Single",
    "Ring Perpendicular to ND",
    "Ring Perpendicular to RD",
    "Ring Perpendicular to TD",
    "OffsetFTD",
    "OffsetATD",
    "OffsetF45",
    "OffsetA45",
    "OffsetFRD",
    "OffsetARD",
    "Gaussian Quadrature",
    "Michael-CCHex",
    "Michael-SNHex"
)

$dataCols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

# Last pre-existing data row (before this edit) - used as a formatting
# template for any brand-new rows appended below it.
$templateRow = 29

for ($i = 0; $i -lt $materials.Length; $i++) {
    $row = $i + 3
    $isNewRow = $row -gt $templateRow

    $ws.Range("A$row").Value = $i + 1
    $ws.Range("B$row").Value = $materials[$i]
    foreach ($col in $dataCols) {
        $ws.Range("$col$row").Value = 1
    }

    if ($isNewRow) {
        # Newly appended rows need column A's bold/bordered/centered style
        # copied over from the existing template row.
        $ws.Range("A$templateRow").Copy()
        $ws.Range("A$row").PasteSpecial(-4122)
    }
}

$excel.CutCopyMode = $false

Write-Output "Rewrote $($materials.Length) material rows (A3:T$($materials.Length + 2))"
